$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "...for internation" + "al business" + " " + "and international
#    representation for local business " + "in a unique exchange..."
#    becomes
#    "...for internation" + "al economic acitivity" + " " +
#    "and international representation for local economic activity " +
#    "in a unique exchange..."
# ------------------------------------------------------------------
$old0 = "al business and international representation for local business "
$new0 = "al economic acitivity and international representation for local economic activity "
$found0 = $d.Content.Find.Execute($old0, $true, $false, $false, $false, $false, $true, 1, $false, $new0, 2)

# ------------------------------------------------------------------
# 2) The trailing part of the sentence ("international interaction
#    and peer-to-peer coordination...wide application. ") is removed
#    from in front of the _GoBack bookmark and replaced with the much
#    shorter "on a global scale as expanding reciprocity" (still in
#    front of the bookmark).
# ------------------------------------------------------------------
$old1 = "international interaction and peer-to-peer coordination…a novel form of grassroots international trade that unites small business while galvanizing consumer participation in an interactive process with broad consideration and wide application. "
$new1 = "on a global scale as expanding reciprocity"
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# ------------------------------------------------------------------
# 3) The text that used to sit in front of the bookmark re-appears
#    right after it (the bookmark itself must not move into/through
#    the new text, so edit the single space run that immediately
#    follows the bookmark instead of spanning across/over it, which
#    would delete the bookmark).
# ------------------------------------------------------------------
$tailText = "peer-to-peer coordination…a novel form of grassroots international trade that unites small business while galvanizing consumer participation in an interactive process with broad consideration and wide application."

$bm = $null
try {
    $bm = $d.Bookmarks("_GoBack")
} catch {
    $bm = $null
}

if ($bm -ne $null) {
    $afterBookmark = $d.Range($bm.Start, $bm.Start + 1)
    $afterBookmark.Text = " is facilitated by " + $tailText + "  "
} else {
    # Fallback: the bookmark was not found by name for some reason -
    # locate the insertion point via the surrounding text instead and
    # simply insert the text there (bookmark preservation best effort).
    $anchor = "on a global scale as expanding reciprocity"
    $full = $d.Content.Text
    $idx = $full.IndexOf($anchor)
    if ($idx -ge 0) {
        $insertionPoint = $d.Range($idx + $anchor.Length, $idx + $anchor.Length)
        $insertionPoint.InsertAfter(" is facilitated by " + $tailText + "  ")
    }
}

$bmFinal = $null
try { $bmFinal = $d.Bookmarks("_GoBack") } catch { $bmFinal = $null }
if ($bmFinal -ne $null) {
    Write-Host "Found0=$found0 Found1=$found1 BookmarkStart=$($bmFinal.Start)"
} else {
    Write-Host "Found0=$found0 Found1=$found1 BookmarkStart=<missing>"
}
